$d = $word.ActiveDocument

# 1. Fix the ">>>" trailing text to ">>" (remove one trailing '>').
$d.Content.Find.Execute("stuff after this line >>>", $true, $false, $false, $false, $false,
                         $true, 1, $false, "stuff after this line >>", 2)

# 2. Insert a new paragraph after the ">>> your ... >>" line with the GitHub comment in red,
#    followed by an empty paragraph.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*stuff after this line*") {
        $target = $p
    }
}

$insertionRange = $target.Range
$insertionRange.Collapse(0)  # wdCollapseEnd

$insertionRange.InsertParagraphAfter()
$insertionRange.InsertParagraphAfter()

# Move to the first newly inserted paragraph
$newPara = $target.Next()
$newPara.Range.Text = "Git hub is so cool. I love being able to manage different version of my software when collaborating with other people regarding personal projects and work.  Hope everyone is doing well with covid and all that Jazz, stay safe people."
$newPara.Range.Font.Color = 255  # wdColorRed (0x0000FF in BGR = red)
